# Update "want to go" counts (column F) across the workbook's sheets.
# Sheet 1 = 展览 (Exhibitions)
# Sheet 2 = 演出 (Performances)
# Sheet 4 = 全部类型 (All types)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value  = 5006
$ws1.Range("F8").Value  = 10
$ws1.Range("F10").Value = 514
$ws1.Range("F13").Value = 1397
$ws1.Range("F14").Value = 3671
$ws1.Range("F17").Value = 121
$ws1.Range("F18").Value = 83
$ws1.Range("F19").Value = 2663
$ws1.Range("F21").Value = 17
$ws1.Range("F27").Value = 60

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 42

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value  = 42
$ws4.Range("F6").Value  = 5006
$ws4.Range("F9").Value  = 10
$ws4.Range("F11").Value = 514
$ws4.Range("F14").Value = 1397
$ws4.Range("F15").Value = 3671
$ws4.Range("F18").Value = 121
$ws4.Range("F19").Value = 83
$ws4.Range("F20").Value = 2663
$ws4.Range("F22").Value = 17
$ws4.Range("F28").Value = 60
